# Apply the taxon-record shuffle described by the diff.
# Columns A,B,D,E,F,G,H hold the per-record taxon data that gets moved
# between rows; columns C and I are constant across these rows and are
# left untouched. Row 2 additionally carries bird-survey columns
# K,L,M,N which must travel together with the rest of its data when it
# swaps with row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    $data = @{}
    foreach ($col in @("A","B","D","E","F","G","H")) {
        $data[$col] = $ws.Range("$col$row").Value2
    }
    return $data
}

function Set-RowData($row, $data) {
    foreach ($col in @("A","B","D","E","F","G","H")) {
        $ws.Range("$col$row").Value = $data[$col]
    }
}

# Capture original values for rows 2-7 before any writes.
$row2 = Get-RowData 2
$row3 = Get-RowData 3
$row4 = Get-RowData 4
$row5 = Get-RowData 5
$row7 = Get-RowData 7

# Also capture the bird-survey columns that live on row 2 originally.
$k2 = $ws.Range("K2").Value2
$l2 = $ws.Range("L2").Value2
$m2 = $ws.Range("M2").Value2
$n2 = $ws.Range("N2").Value2

# Row 2 <-> Row 7 full swap (taxon data + the K/L/M/N activity columns).
Set-RowData 2 $row7
Set-RowData 7 $row2

$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""

$ws.Range("K7").Value = $k2
$ws.Range("L7").Value = $l2
$ws.Range("M7").Value = $m2
$ws.Range("N7").Value = $n2

# Rows 3,4,5 cyclic rotation: 3<-4, 4<-5, 5<-3.
Set-RowData 3 $row4
Set-RowData 4 $row5
Set-RowData 5 $row3

# Row 6 is unchanged.
